$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New workshop entry (row 37) - "Time-to-Event Analysis" (Spring 2022)
$ws.Range("A37").Value = 6
$ws.Range("B37").Value = "set"
$ws.Range("E37").Value = "Spring 2022"
$ws.Range("C37").Value = "Time-to-Event Analysis, aka Survival Analysis"
$ws.Range("D37").Value = "Seven 1-hour sessions"
$ws.Range("G37").Value = "Inferior Models: You must make a choice - censoring and duration of follow-up"
$ws.Range("G38").Value = "Non-parametric approach: Kaplan Meyer Plots with Log Rank Test (bivariate)"
$ws.Range("G39").Value = "Semi-parametric approach: Cox Proportional-Hazards Regression (multivariate)"
$ws.Range("G40").Value = "Extension: Time varying covariates in a Cox model"

# New small bold style used for the sub-bullet rows (G38:G40)
$st = $wb.Styles.Add("WorkshopSubStyle")
$st.Font.Bold = $true
$st.Font.Size = 6
$st.Font.Color = 2696481
$st.Font.Name = "Segoe UI"
$ws.Range("G38:G40").Style = "WorkshopSubStyle"
$st.Delete()

# Match the recorded UI state: selection on the newly added cell
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G40").Select()

# Page setup was touched (printed/previewed) in the source edit
$ws.PageSetup.Orientation = 1
